$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-13
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03)
$ws.Range("C2:C13").Value = 45233
